$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (pushing existing rows 35.. down by one)
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with values (mirrors the row directly
# below it, except for the new date in column D)
$ws.Range("A35").Value = 3
$ws.Range("B35").Value = "Femacal de La Calera"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44547
$ws.Range("D35").NumberFormat = $ws.Range("D36").NumberFormat
$ws.Range("E35").Value = 5
$ws.Range("F35").Value = 100112039
$ws.Range("G35").Value = "Ciboulette"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = 1500
$ws.Range("N35").Value = "`$/docena de atados"
$ws.Range("O35").Value = "Provincia de Quillota"
$ws.Range("P35").Value = 500
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = "Hortaliza"
